$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Item "DRAMENEX 50MG 20 TABS." (row 25) transaction ratio corrected from 2:1 to 1:1
$ws.Cells.Item(25, 8).Value = "1:1"

# Remove the now-stale stock rows (items 49-74, rows 55-80) together with the
# subtotal row (81) that summed their sale price column. Deleting the whole
# block shifts the trailing footer (timestamp / page / credit line), that used
# to live on row 82, up to row 55.
$ws.Range("A55:Q81").EntireRow.Delete()

# Refresh the export timestamp shown in the footer to the new save time.
$ws.Cells.Item(55, 1).Value = "Sunday, 22 June, 2025 6:29 PM"
